# "added july 20 data"
# Update the team stat figures for the four teams whose columns changed
# (C = Los Angeles Lakers, E = Los Angeles Clippers, I = Utah Jazz,
#  S = New Orleans Pelicans) across every stat row (2-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 56.1
$ws.Range("E2").Value = 51.4
$ws.Range("I2").Value = 59.5
$ws.Range("S2").Value = 67

$ws.Range("C3").Value = 43.9
$ws.Range("E3").Value = 48.6
$ws.Range("I3").Value = 40.5
$ws.Range("S3").Value = 33

$ws.Range("C4").Value = 40.8
$ws.Range("E4").Value = 31.7
$ws.Range("I4").Value = 54.7
$ws.Range("S4").Value = 55.8

$ws.Range("C5").Value = 1.9
$ws.Range("E5").Value = 7.9
$ws.Range("I5").Value = 1.9
$ws.Range("S5").Value = 9.6

$ws.Range("C6").Value = 32
$ws.Range("E6").Value = 47.5
$ws.Range("I6").Value = 22.6
$ws.Range("S6").Value = 31.7

$ws.Range("C7").Value = 8.7
$ws.Range("E7").Value = 7.9
$ws.Range("I7").Value = 6.6
$ws.Range("S7").Value = 6.7

$ws.Range("C8").Value = 27.2
$ws.Range("E8").Value = 20.8
$ws.Range("I8").Value = 22.6
$ws.Range("S8").Value = 12.5

$ws.Range("C9").Value = 13.6
$ws.Range("E9").Value = 28.7
$ws.Range("I9").Value = 18.9
$ws.Range("S9").Value = 13.5

$ws.Range("C10").Value = 38.8
$ws.Range("E10").Value = 23.8
$ws.Range("I10").Value = 52.8
$ws.Range("S10").Value = 46.2

$ws.Range("C11").Value = 52.4
$ws.Range("E11").Value = 50
$ws.Range("I11").Value = 41.4
$ws.Range("S11").Value = 41.4

$ws.Range("C12").Value = 47.6
$ws.Range("E12").Value = 50
$ws.Range("I12").Value = 58.6
$ws.Range("S12").Value = 58.6

$ws.Range("C13").Value = 90.9
$ws.Range("E13").Value = 56.3
$ws.Range("I13").Value = 62.5
$ws.Range("S13").Value = 90.9

$ws.Range("C14").Value = 9.1
$ws.Range("E14").Value = 43.8
$ws.Range("I14").Value = 37.5
$ws.Range("S14").Value = 9.1

$ws.Range("C15").Value = 65.6
$ws.Range("E15").Value = 53.1
$ws.Range("I15").Value = 45.9
$ws.Range("S15").Value = 55

$ws.Range("C16").Value = 34.4
$ws.Range("E16").Value = 46.9
$ws.Range("I16").Value = 54.1
$ws.Range("S16").Value = 45

# The leftover placeholder rows below the table (rows 19-31) are tidied up:
# rows 19-23 are fully cleared (no longer used at all), and rows 24-31 have
# their now-unused leading TEAM/value columns (A-D) cleared while the
# remaining formatted columns (E-P) are left in place.
$ws.Range("A19:P23").Clear()
$ws.Range("A24:D31").Clear()

# Leave the selection where the user ended up after this edit.
$ws.Range("F30").Select() | Out-Null
